# ---------------------------------------------------------------------------
# Commit: "updated the sort by difficulty"
#
# 1. On the "Recipes" sheet, the Cuisine (C) / Skill (D) / Time (E) columns
#    for rows 8, 9 and 10 had been entered one column too far to the right
#    (e.g. the Skill value was sitting in column C, Time in D, Cuisine in
#    E). This rotates each row's C/D/E values one step to the left so the
#    new C,D,E = old D, old E, old C.
# 2. A stray double-space in row 8's "Alt" column is collapsed to one space.
# 3. Row 9's Ingredient column is rewritten from a plain comma separated
#    string into a python-list "repr" style string.
# 4. Three new per-recipe ingredient-breakdown sheets are appended, named
#    "10", "9" and "8" (matching recipe IDs 10, 9 and 8 respectively), each
#    with the same Quantity / Measurement / Ingredient layout used by the
#    other per-recipe sheets already in the workbook.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$recipes = $wb.Worksheets.Item("Recipes")

# --- 1 & 2 & 3: fix up the Recipes sheet ------------------------------------

# Row 8 - Chicken Wraps (recipe id 8): rotate Cuisine/Skill/Time left
$recipes.Range("C8").Value = "Intermediate"
$recipes.Range("D8").Value = "30 Mins"
$recipes.Range("E8").Value = "Western"
$recipes.Range("H8").Value = "Bacon --> Ham, Ranch Dressing --> Mayonnaise"

# Row 9 - Tom Yum Soup (recipe id 9): rotate Cuisine/Skill/Time left
$recipes.Range("C9").Value = "Hard"
$recipes.Range("D9").Value = "1.5 Hours"
$recipes.Range("E9").Value = "Thai"
$recipes.Range("G9").Value = "['240 ml Chicken Stock', '0.5 Stalk Lemongrass (Halved)', '0.75 Galangal (Sliced)', '10 g Tom Yum Chilli Paste', '0.5 Kaffir Lime Leaf', '3 Oyster Mushrooms', '0.25 Thai Chilli Pepper (Halved)', '2 Large Shrimp', '3 g Sugar', '4 ml Lime Juice', '10 g Fresh Cilantro Leaves']"

# Row 10 - Steamed Eggs (recipe id 10): rotate Cuisine/Skill/Time left
$recipes.Range("C10").Value = "Intermediate"
$recipes.Range("D10").Value = "45 Mins"
$recipes.Range("E10").Value = "Korean"

# --- 4: append the three new per-recipe ingredient sheets -------------------

function Add-IngredientSheet {
    param($afterSheet, $sheetName, $rows)

    $newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
    $newSheet.Name = $sheetName

    $newSheet.Cells.Item(1, 1).Value = "Quantity"
    $newSheet.Cells.Item(1, 2).Value = "Measurement"
    $newSheet.Cells.Item(1, 3).Value = "Ingredient"
    $newSheet.Range("A1:C1").Font.Bold = $true

    $r = 2
    foreach ($row in $rows) {
        $newSheet.Cells.Item($r, 1).Value = $row[0]
        $newSheet.Cells.Item($r, 2).Value = $row[1]
        $newSheet.Cells.Item($r, 3).Value = $row[2]
        $r = $r + 1
    }

    return $newSheet
}

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# New sheet "10" -> recipe id 10 (Steamed Eggs) ingredient breakdown
$steamedEggsRows = @(
    , @(2,    "N/A", "Large Eggs")
    , @(240,  "ml",  "Chicken Broth")
    , @(1,    "N/A", "Tablespoon Light Soy Sauce")
    , @(1,    "N/A", "Teaspoon Sesame Oil")
    , @(0.25, "N/A", "Teaspoon Salt")
    , @(10,   "g",   "Green Onions")
    , @(0.25, "N/A", "Teaspoon White Pepper")
)
$sheet10tab = Add-IngredientSheet $lastSheet "10" $steamedEggsRows

# New sheet "9" -> recipe id 9 (Tom Yum Soup) ingredient breakdown
$tomYumRows = @(
    , @(240,  "ml",  "Chicken Stock'")
    , @(0.5,  "N/A", "Stalk Lemongrass")
    , @(0.75, "N/A", "Galangal")
    , @(10,   "g",   "Tom Yum Chilli Paste'")
    , @(0.5,  "N/A", "Kaffir Lime Leaf'")
    , @(3,    "N/A", "Oyster Mushrooms'")
    , @(0.25, "N/A", "Thai Chilli Pepper")
    , @(2,    "N/A", "Large Shrimp'")
    , @(3,    "g",   "Sugar'")
    , @(4,    "ml",  "Lime Juice'")
    , @(10,   "g",   "Fresh Cilantro Leaves']")
)
$sheet9tab = Add-IngredientSheet $sheet10tab "9" $tomYumRows

# New sheet "8" -> recipe id 8 (Chicken Wraps) ingredient breakdown
$chickenWrapRows = @(
    , @(1,    "N/A", "Tortilla")
    , @(20,   "g",   "Chicken Breast")
    , @(30,   "g",   "Lettuce")
    , @(2,    "N/A", "Slices of Bacon")
    , @(0.25, "N/A", "Tomato")
    , @(0.5,  "N/A", "Onion")
    , @(0.25, "N/A", "Avocado")
    , @(1.5,  "N/A", "Tablespoon Shredded Cheese")
    , @(30,   "g",   "Ranch Dressing")
    , @(1,    "N/A", "Teaspoon Hot Sauce")
)
$sheet8tab = Add-IngredientSheet $sheet9tab "8" $chickenWrapRows

# Leave selection back on the Recipes sheet, matching the original file.
$recipes.Activate()
